# Update "想去人数" (F column) figures on the "展览" and "全部类型" sheets
# to reflect the latest scrape (gh-pages output generated at 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F2").Value = 5443
$wsExhibit.Range("F6").Value = 822
$wsExhibit.Range("F7").Value = 21
$wsExhibit.Range("F8").Value = 348

# --- Sheet "全部类型" ---
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 5443
$wsAll.Range("F6").Value = 822
$wsAll.Range("F7").Value = 21
$wsAll.Range("F9").Value = 348
